# Refactored data preparation pipeline:
# Swap the "promotion"/"distribution" columns (D <-> E) and the
# "off_trade_visibility"/"covid" columns (G <-> H), for the header
# row and every data row on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 158 }

for ($r = 1; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $gCell = $ws.Cells.Item($r, 7)
    $hCell = $ws.Cells.Item($r, 8)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $gVal = $gCell.Value2
    $hVal = $hCell.Value2

    $dCell.Value2 = $eVal
    $eCell.Value2 = $dVal

    $gCell.Value2 = $hVal
    $hCell.Value2 = $gVal
}
